# Add an "Email" column (H) to the suppliers sheet: header in H1, the
# supplier's email address in H2 — mirrors the existing A:G header/value
# layout used for Name/Registration Number/CIF/Address/Bank/IBAN/amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Email"
$ws.Range("H2").Value = "adrianrentea01@gmail.com"

# Size the new column to fit its contents, like the other bestFit columns.
$ws.Columns.Item(8).ColumnWidth = 23.6

# Leave the new cells selected, as a user adding this column interactively
# would after typing the two values in.
$ws.Range("H1:H2").Select()
